# Fruta / hortaliza, semanal
# Insert a new weekly record row (row 110) into the Cilantro sheet, pushing
# the existing rows 110..202 down to 111..203.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110 (shifts rows 110-202 down to 111-203,
# and copies formatting - including the date style on column D - from the
# row above, same as native Excel "Insert Copied Cells"/"Insert Sheet Rows").
$ws.Rows.Item(110).Insert()

# Populate the new row 110 with the new weekly data point.
$ws.Range("A110").Value = 8
$ws.Range("B110").Value = "Terminal La Palmera de La Serena"
$ws.Range("C110").Value = "Coquimbo"
$ws.Range("D110").Value = 44957
$ws.Range("E110").Value = 4
$ws.Range("F110").Value = 100112040
$ws.Range("G110").Value = "Cilantro"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = 3250
$ws.Range("N110").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O110").Value = "Provincia del Elquí"
$ws.Range("P110").Value = 2167
$ws.Range("Q110").Value = 1.5
$ws.Range("R110").Value = "Hortaliza"
